$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D8").Value = "Hình thức:"
$ws.Range("D9").Value = "Công khai:"
$ws.Range("D10").Value = "Giai đoạn:"

$ws.Range("D10").Select()
